# Applies the "Obsidian UDSE Plan Phase 01" update:
#   1. Refresh the cached "datetimeFigureOut" footer date (15/05/2020 -> 09/10/2020)
#      on the slide master and every slide layout.
#   2. Rename the "Relationships" keyword box to "Abstractions".
#   3. Rename the "Causal Chains" keyword box to "Causal Relationships"
#      (kept as two runs: "Causal " + "Relationships").

$p = $ppt.ActivePresentation

$ppPlaceholderDate = 16
$newDate = "09/10/2020"

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            $isDatePh = $false
            if ($shp.Type -eq 14) {
                if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                    $isDatePh = $true
                }
            }
            if ($isDatePh) {
                if ($shp.TextFrame.TextRange.Text -ne $newDate) {
                    $shp.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

# 1a. Slide master footer date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# 1b. Every slide layout's footer date placeholder.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# 2 & 3. Update the keyword text boxes on slide 1.
$slide = $p.Slides.Item(1)

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $txt = $shp.TextFrame.TextRange.Text
        if ($txt -eq "Relationships") {
            $shp.TextFrame.TextRange.Text = "Abstractions"
        }
        elseif ($txt -eq "Causal Chains") {
            $tr = $shp.TextFrame.TextRange
            $tr.Text = "Causal Relationships"
            # Split into two runs ("Causal " + "Relationships") matching the
            # authored edit, while keeping identical run formatting.
            $firstPart = $tr.Characters(1, 7)
            $firstPart.Text = "Causal "
        }
    }
}
